$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common values shared across all data rows (2-13)
$projectId = "PRJ4897790"
$projectName = "2025 Planning and Scheduling IDT Work"
$ratingDate = "2025-09-23"

for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = $projectId
    $ws.Cells.Item($r, 2).Value = $projectName
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $ratingDate
}

# Row 2 - firmness_of_opportunity
$ws.Cells.Item(2, 5).Value = "Red"
$ws.Cells.Item(2, 6).Value = "DAF / GIP or 'one pager' not yet prepared. Approvals not expected for >3 months. Start date of IT effort unknown."

# Row 3 - business_case
$ws.Cells.Item(3, 5).Value = "Green"
$ws.Cells.Item(3, 6).Value = "The project status report includes benefits described as 'The 'HVO B7' tool is projected to generate approximately `$3 million per annum in tax rebates for Shell, with the value poised to grow as sales of HVO B7 continue to increase.' Additionally, '31 m`$ / year level 1/level 2 Europe benefits were confirmed, down from 55 m`$/year.' This indicates the presence of L1 and/or L2 benefits as defined by the project with agreed benefit ownership."

# Row 4 - business_requirements
$ws.Cells.Item(4, 5).Value = "Green"
$ws.Cells.Item(4, 6).Value = "The project closure activities are in progress with successful completion of Hypercare and business sign-off, reflecting well-documented requirements for impacted portfolios. The consistent and ongoing completion of tasks and achievements, such as finalized technical documentation and agreed deployment plans, supports readiness to mobilize."

# Row 5 - solution
$ws.Cells.Item(5, 5).Value = "Green"
$ws.Cells.Item(5, 6).Value = "Progress & Success: Hypercare successfully completed and ATO signed-off as per plan. Project team dispersed. Closure activities are in progress. The 'HVO B7' tool is projected to generate approximately `$3 million per annum in tax rebates for Shell, with the value poised to grow as sales of HVO B7 continue to increase. Challenges & RTGs: None"

# Row 6 - estimate_budget_accuracy
$ws.Cells.Item(6, 5).Value = "Red"
$ws.Cells.Item(6, 6).Value = "Cost (RED): `$750k budget approved through to the end of July has been spent. RTG: Sandeep to request bridge funding of `$250k cover until end of Sep'25 and GIP approval."

# Row 7 - commercialops_functions_support (rating unchanged: Green)
$ws.Cells.Item(7, 6).Value = "Project team already provided hand over to the support team regarding the changes they need to make. This implies that existing staff can deliver all necessary work as the project is in its closing phase and no specific skill gaps or resource needs are mentioned."

# Row 8 - ts_and_broader_rds_it_resource_Capacity
$ws.Cells.Item(8, 5).Value = "Green"
$ws.Cells.Item(8, 6).Value = "Worked with CC team and managed to retain the primary data engineer to continue work on the project, hence the associated risk is mitigated now and Resourcing status moved back to Green."

# Row 9 - ts_process_team_within_comm_ops_bia_org (rating unchanged: Green)
$ws.Cells.Item(9, 6).Value = "Project team already provided hand over to the support team regarding the changes they need to make. This implies that existing staff can deliver all necessary work as the project is in its closing phase and no specific skill gaps or resource needs are mentioned."

# Row 10 - owning_portfolio_it_resource_capacity (rating unchanged: Green)
$ws.Cells.Item(10, 6).Value = "Project team already provided hand over to the support team regarding the changes they need to make. This implies that existing staff can deliver all necessary work as the project is in its closing phase and no specific resource shortages or skill gaps are mentioned."

# Row 11 - ongoing_supportability_of_solution
$ws.Cells.Item(11, 5).Value = "Green"
$ws.Cells.Item(11, 6).Value = "Project team already provided hand over to the support team regarding the changes they need to make. Additionally, Hypercare was successfully completed with sign-off from Business, which indicates that the solution fits within an existing support model or contract, fulfilling the criteria for Green."

# Row 12 - effective_governance (rating unchanged: Red)
$ws.Cells.Item(12, 6).Value = "The project status report does not contain any mention of a governance model being defined, drafted, or agreed upon, nor a clear decision owner identified. Therefore, the rating is Red because 'Governance model not yet drafted'."

# Row 13 - partnership_and_collaboration
$ws.Cells.Item(13, 5).Value = "Green"
$ws.Cells.Item(13, 6).Value = "Hypercare successfully completed and received sign-off from Business. Additionally, there is evidence of coordinated efforts between various teams, such as the AMDP team working with T&S support to make appropriate changes, and technical documentation being shared with the support team. These imply a strong working relationship between business and IT teams, meeting the criteria for Green."
